{"js": "// Update each two-digit multiplication equation cell to its new value.\n// Every old equation string is unique in the document, so a direct\n// search-and-replace (matchCase, no wildcards) for each pair is safe.\nconst replacements = [\n  [\"55\u00d747=2585\", \"77\u00d750=3850\"],\n  [\"37\u00d711=407\", \"37\u00d716=592\"],\n  [\"86\u00d749=4214\", \"30\u00d723=690\"],\n  [\"81\u00d738=3078\", \"89\u00d796=8544\"],\n  [\"70\u00d720=1400\", \"55\u00d717=935\"],\n  [\"54\u00d712=648\", \"62\u00d741=2542\"],\n  [\"36\u00d741=1476\", \"89\u00d783=7387\"],\n  [\"11\u00d752=572\", \"26\u00d739=1014\"],\n  [\"89\u00d782=7298\", \"30\u00d712=360\"],\n  [\"60\u00d768=4080\", \"12\u00d771=852\"],\n  [\"13\u00d719=247\", \"20\u00d789=1780\"],\n  [\"87\u00d786=7482\", \"25\u00d760=1500\"],\n  [\"48\u00d799=4752\", \"43\u00d720=860\"],\n  [\"44\u00d768=2992\", \"13\u00d794=1222\"],\n  [\"69\u00d736=2484\", \"97\u00d749=4753\"],\n  [\"78\u00d771=5538\", \"41\u00d770=2870\"],\n  [\"71\u00d729=2059\", \"81\u00d799=8019\"],\n  [\"23\u00d750=1150\", \"14\u00d726=364\"],\n  [\"32\u00d726=832\", \"21\u00d745=945\"],\n  [\"58\u00d779=4582\", \"40\u00d797=3880\"],\n  [\"49\u00d743=2107\", \"35\u00d752=1820\"],\n  [\"60\u00d748=2880\", \"47\u00d785=3995\"],\n  [\"27\u00d778=2106\", \"79\u00d755=4345\"],\n  [\"77\u00d752=4004\", \"29\u00d742=1218\"],\n  [\"22\u00d723=506\", \"91\u00d732=2912\"]\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update each two-digit multiplication equation cell to its new value.\n# Every old equation string is unique in the document, so a direct\n# Find/Replace (MatchCase, no wildcards) for each pair is safe and\n# leaves run formatting (font/size) untouched.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"55\u00d747=2585\"; New = \"77\u00d750=3850\" },\n    @{ Old = \"37\u00d711=407\"; New = \"37\u00d716=592\" },\n    @{ Old = \"86\u00d749=4214\"; New = \"30\u00d723=690\" },\n    @{ Old = \"81\u00d738=3078\"; New = \"89\u00d796=8544\" },\n    @{ Old = \"70\u00d720=1400\"; New = \"55\u00d717=935\" },\n    @{ Old = \"54\u00d712=648\"; New = \"62\u00d741=2542\" },\n    @{ Old = \"36\u00d741=1476\"; New = \"89\u00d783=7387\" },\n    @{ Old = \"11\u00d752=572\"; New = \"26\u00d739=1014\" },\n    @{ Old = \"89\u00d782=7298\"; New = \"30\u00d712=360\" },\n    @{ Old = \"60\u00d768=4080\"; New = \"12\u00d771=852\" },\n    @{ Old = \"13\u00d719=247\"; New = \"20\u00d789=1780\" },\n    @{ Old = \"87\u00d786=7482\"; New = \"25\u00d760=1500\" },\n    @{ Old = \"48\u00d799=4752\"; New = \"43\u00d720=860\" },\n    @{ Old = \"44\u00d768=2992\"; New = \"13\u00d794=1222\" },\n    @{ Old = \"69\u00d736=2484\"; New = \"97\u00d749=4753\" },\n    @{ Old = \"78\u00d771=5538\"; New = \"41\u00d770=2870\" },\n    @{ Old = \"71\u00d729=2059\"; New = \"81\u00d799=8019\" },\n    @{ Old = \"23\u00d750=1150\"; New = \"14\u00d726=364\" },\n    @{ Old = \"32\u00d726=832\"; New = \"21\u00d745=945\" },\n    @{ Old = \"58\u00d779=4582\"; New = \"40\u00d797=3880\" },\n    @{ Old = \"49\u00d743=2107\"; New = \"35\u00d752=1820\" },\n    @{ Old = \"60\u00d748=2880\"; New = \"47\u00d785=3995\" },\n    @{ Old = \"27\u00d778=2106\"; New = \"79\u00d755=4345\" },\n    @{ Old = \"77\u00d752=4004\"; New = \"29\u00d742=1218\" },\n    @{ Old = \"22\u00d723=506\"; New = \"91\u00d732=2912\" }\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $pair.Old\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $pair.New\n\n    $found = $find.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $null, 2)\n    if (-not $found) {\n        throw \"No match found for: $($pair.Old)\"\n    }\n}\n\nWrite-Output \"done\"\n"}
